$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status text for row 2 and row 3 ("OPTIMAL" -> "TIME_LIMIT")
$ws.Range("E2").Value = "TIME_LIMIT"
$ws.Range("E3").Value = "TIME_LIMIT"

# Correct the fixed recourse data values for rows 2 and 3
$ws.Range("B2").Value = -636.8209940539369
$ws.Range("C2").Value = 6.8750427571620385
$ws.Range("D2").Value = 3623.924046705

$ws.Range("B3").Value = -642.701897941751
$ws.Range("C3").Value = 9.5353207677372
$ws.Range("D3").Value = 3766.74439455

# Remove rows 4 through 11 entirely, shrinking the used range to A1:H3
$ws.Range("A4:H11").Delete()
